# Apply updated cryptos data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '60.255.25'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '2.589.99'
$ws.Range("E3").Value = '  -1.88%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''508.61'
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").Value = '''153.75'
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("D7").Value = '''1.01'
$ws.Range("E7").Value = '  +1.28%  '
$ws.Range("E8").Value = '  -3.61%  '
$ws.Range("D9").Value = '2.600.13'
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").Value = '''6.72'
$ws.Range("E10").Value = '  +7.11%  '
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '3.046.46'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '60.251.04'
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '2.596.47'
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("E19").Value = '  -1.41%  '
$ws.Range("D20").Value = '''352.54'
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").Value = '''10.53'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '''6.11'
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").Value = '''60.37'
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = '''0.420'
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").Value = '''0.166'
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("D28").Value = '0.0₃0837'
$ws.Range("E28").Value = '  -3.99%  '
$ws.Range("D29").Value = '''7.33'
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").Value = '''151.83'
$ws.Range("E32").Value = '  -3.41%  '
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("D35").Value = '''3.99'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -3.38%  '
$ws.Range("D37").Value = '''0.865'
$ws.Range("E37").Value = '  +4.32%  '
$ws.Range("E38").Value = '  -3.54%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '''36.09'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '''0.838'
$ws.Range("E40").Value = '  -3.97%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '''3.75'
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("D42").Value = '''296.05'
$ws.Range("E42").Value = '  -4.51%  '
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("D46").Value = '''0.0552'
$ws.Range("E46").Value = '  -4.82%  '
$ws.Range("D47").Value = '''19.71'
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").Value = '''4.78'
$ws.Range("E48").Value = '  -4.80%  '
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").Value = '''10.31'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = '1.990.25'
$ws.Range("E51").Value = '  -2.88%  '
